$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# --- Sheet "Recommandations": update rows 2-29, then delete old rows 30-34 ---
$ws1.Range("A2").Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws1.Range("B2").Value = 0
$ws1.Range("C2").Value = 4
$ws1.Range("D2").Value = 647.41
$ws1.Range("E2").Value = 164.65
$ws1.Range("F2").Value = "🟡 Observer"
$ws1.Range("G2").Value = "➖ Neutre"

$ws1.Range("A3").Value = "BRVM - SERVICES FINANCIERS"
$ws1.Range("B3").Value = 0
$ws1.Range("C3").Value = 4
$ws1.Range("D3").Value = 583.13
$ws1.Range("E3").Value = 145.01
$ws1.Range("F3").Value = "🟡 Observer"
$ws1.Range("G3").Value = "➖ Neutre"

$ws1.Range("A4").Value = "BRVM-PRESTIGE"
$ws1.Range("B4").Value = 0
$ws1.Range("C4").Value = 4
$ws1.Range("D4").Value = 571.37
$ws1.Range("E4").Value = 142.04
$ws1.Range("F4").Value = "🟡 Observer"
$ws1.Range("G4").Value = "➖ Neutre"

$ws1.Range("A5").Value = "BRVM - INDUSTRIELS"
$ws1.Range("B5").Value = 0
$ws1.Range("C5").Value = 4
$ws1.Range("D5").Value = 535.5
$ws1.Range("E5").Value = 136.25
$ws1.Range("F5").Value = "🟡 Observer"
$ws1.Range("G5").Value = "➖ Neutre"

$ws1.Range("A6").Value = "BRVM - ENERGIE"
$ws1.Range("B6").Value = 0
$ws1.Range("C6").Value = 4
$ws1.Range("D6").Value = 447.65
$ws1.Range("E6").Value = 112.34
$ws1.Range("F6").Value = "🟡 Observer"
$ws1.Range("G6").Value = "➖ Neutre"

$ws1.Range("A7").Value = "BRVM - SERVICES PUBLICS"
$ws1.Range("B7").Value = 0
$ws1.Range("C7").Value = 4
$ws1.Range("D7").Value = 428.4
$ws1.Range("E7").Value = 107.16
$ws1.Range("F7").Value = "🟡 Observer"
$ws1.Range("G7").Value = "➖ Neutre"

$ws1.Range("A8").Value = "BRVM - TELECOMMUNICATIONS"
$ws1.Range("B8").Value = 0
$ws1.Range("C8").Value = 4
$ws1.Range("D8").Value = 373.56
$ws1.Range("E8").Value = 93.28
$ws1.Range("F8").Value = "🟡 Observer"
$ws1.Range("G8").Value = "➖ Neutre"

$ws1.Range("A9").Value = "EVIOSYS PACKAGING SIEM CI (SEMC)"
$ws1.Range("B9").Value = 4
$ws1.Range("C9").Value = 0
$ws1.Range("D9").Value = 28.49
$ws1.Range("E9").Value = 7.07
$ws1.Range("F9").Value = "🟢 Achat"
$ws1.Range("G9").Value = "✅ Renforcer"

$ws1.Range("A10").Value = "SICOR CI (SICC)"
$ws1.Range("B10").Value = 2
$ws1.Range("C10").Value = 0
$ws1.Range("D10").Value = 11.48
$ws1.Range("E10").Value = 4.08
$ws1.Range("F10").Value = "🟡 Observer"
$ws1.Range("G10").Value = "➖ Neutre"

$ws1.Range("A11").Value = "NEI-CEDA CI (NEIC)"
$ws1.Range("B11").Value = 2
$ws1.Range("C11").Value = 1
$ws1.Range("D11").Value = 11.03
$ws1.Range("E11").Value = 6.48
$ws1.Range("F11").Value = "🟡 Observer"
$ws1.Range("G11").Value = "👀 À surveiller"

$ws1.Range("A12").Value = "VIVO ENERGY CI (SHEC)"
$ws1.Range("B12").Value = 2
$ws1.Range("C12").Value = 0
$ws1.Range("D12").Value = 9.69
$ws1.Range("E12").Value = 4.73
$ws1.Range("F12").Value = "🟡 Observer"
$ws1.Range("G12").Value = "➖ Neutre"

$ws1.Range("A13").Value = "FILTISAC CI (FTSC)"
$ws1.Range("B13").Value = 2
$ws1.Range("C13").Value = 1
$ws1.Range("D13").Value = 7.08
$ws1.Range("E13").Value = 7.26
$ws1.Range("F13").Value = "🟡 Observer"
$ws1.Range("G13").Value = "👀 À surveiller"

$ws1.Range("A14").Value = "SOGB CI (SOGC)"
$ws1.Range("B14").Value = 1
$ws1.Range("C14").Value = 0
$ws1.Range("D14").Value = 3.24
$ws1.Range("E14").Value = 3.24
$ws1.Range("F14").Value = "🟡 Observer"
$ws1.Range("G14").Value = "➖ Neutre"

$ws1.Range("A15").Value = "TOTALENERGIES MARKETING SN (TTLS)"
$ws1.Range("B15").Value = 1
$ws1.Range("C15").Value = 0
$ws1.Range("D15").Value = 3.19
$ws1.Range("E15").Value = 3.19
$ws1.Range("F15").Value = "🟡 Observer"
$ws1.Range("G15").Value = "➖ Neutre"

$ws1.Range("A16").Value = "SICABLE CI (CABC)"
$ws1.Range("B16").Value = 1
$ws1.Range("C16").Value = 1
$ws1.Range("D16").Value = 1.48
$ws1.Range("E16").Value = 3.82
$ws1.Range("F16").Value = "🟡 Observer"
$ws1.Range("G16").Value = "👀 À surveiller"

$ws1.Range("A17").Value = "SOCIETE IVOIRIENNE DE BANQUE  (SIBC)"
$ws1.Range("B17").Value = 1
$ws1.Range("C17").Value = 1
$ws1.Range("D17").Value = 1.42
$ws1.Range("E17").Value = -2.68
$ws1.Range("F17").Value = "🟡 Observer"
$ws1.Range("G17").Value = "👀 À surveiller"

$ws1.Range("A18").Value = "CFAO MOTORS CI (CFAC)"
$ws1.Range("B18").Value = 1
$ws1.Range("C18").Value = 1
$ws1.Range("D18").Value = 0.5
$ws1.Range("E18").Value = 4.71
$ws1.Range("F18").Value = "🟡 Observer"
$ws1.Range("G18").Value = "👀 À surveiller"

$ws1.Range("A19").Value = "UNILEVER CI (UNLC)"
$ws1.Range("B19").Value = 1
$ws1.Range("C19").Value = 1
$ws1.Range("D19").Value = 0.47
$ws1.Range("E19").Value = -6.67
$ws1.Range("F19").Value = "🟡 Observer"
$ws1.Range("G19").Value = "👀 À surveiller"

$ws1.Range("A20").Value = "TOTALENERGIES MARKETING CI (TTLC)"
$ws1.Range("B20").Value = 0
$ws1.Range("C20").Value = 1
$ws1.Range("D20").Value = -1.49
$ws1.Range("E20").Value = -1.49
$ws1.Range("F20").Value = "🟡 Observer"
$ws1.Range("G20").Value = "➖ Neutre"

$ws1.Range("A21").Value = "SERVAIR ABIDJAN CI (ABJC)"
$ws1.Range("B21").Value = 0
$ws1.Range("C21").Value = 1
$ws1.Range("D21").Value = -1.8
$ws1.Range("E21").Value = -1.8
$ws1.Range("F21").Value = "🟡 Observer"
$ws1.Range("G21").Value = "➖ Neutre"

$ws1.Range("A22").Value = "ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)"
$ws1.Range("B22").Value = 1
$ws1.Range("C22").Value = 2
$ws1.Range("D22").Value = -1.92
$ws1.Range("E22").Value = 3.75
$ws1.Range("F22").Value = "🟡 Observer"
$ws1.Range("G22").Value = "👀 À surveiller"

$ws1.Range("A23").Value = "UNIWAX CI (UNXC)"
$ws1.Range("B23").Value = 0
$ws1.Range("C23").Value = 1
$ws1.Range("D23").Value = -2.13
$ws1.Range("E23").Value = -2.13
$ws1.Range("F23").Value = "🟡 Observer"
$ws1.Range("G23").Value = "➖ Neutre"

$ws1.Range("A24").Value = "LOTERIE NATIONALE DU BENIN (LNBB)"
$ws1.Range("B24").Value = 0
$ws1.Range("C24").Value = 1
$ws1.Range("D24").Value = -2.21
$ws1.Range("E24").Value = -2.21
$ws1.Range("F24").Value = "🟡 Observer"
$ws1.Range("G24").Value = "➖ Neutre"

$ws1.Range("A25").Value = "SMB CI (SMBC)"
$ws1.Range("B25").Value = 0
$ws1.Range("C25").Value = 1
$ws1.Range("D25").Value = -3.16
$ws1.Range("E25").Value = -3.16
$ws1.Range("F25").Value = "🟡 Observer"
$ws1.Range("G25").Value = "➖ Neutre"

$ws1.Range("A26").Value = "SOLIBRA CI (SLBC)"
$ws1.Range("B26").Value = 1
$ws1.Range("C26").Value = 2
$ws1.Range("D26").Value = -3.92
$ws1.Range("E26").Value = -4.33
$ws1.Range("F26").Value = "🟡 Observer"
$ws1.Range("G26").Value = "👀 À surveiller"

$ws1.Range("A27").Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws1.Range("B27").Value = 0
$ws1.Range("C27").Value = 1
$ws1.Range("D27").Value = -4.35
$ws1.Range("E27").Value = -4.35
$ws1.Range("F27").Value = "🟡 Observer"
$ws1.Range("G27").Value = "➖ Neutre"

$ws1.Range("A28").Value = "SOCIETE GENERALE COTE D'IVOIRE (SGBC)"
$ws1.Range("B28").Value = 0
$ws1.Range("C28").Value = 2
$ws1.Range("D28").Value = -6.57
$ws1.Range("E28").Value = -3.51
$ws1.Range("F28").Value = "🟡 Observer"
$ws1.Range("G28").Value = "➖ Neutre"

$ws1.Range("A29").Value = "SETAO CI (STAC)"
$ws1.Range("B29").Value = 0
$ws1.Range("C29").Value = 2
$ws1.Range("D29").Value = -12.54
$ws1.Range("E29").Value = -7.05
$ws1.Range("F29").Value = "🟡 Observer"
$ws1.Range("G29").Value = "➖ Neutre"

# Remove now-obsolete trailing rows (30-34) so dimension shrinks to A1:G29
$ws1.Rows("30:34").Delete()

# --- Sheet "Top_YTD": update rows 2-11 ---
$ws2.Range("A2").Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws2.Range("B2").Value = 4600.58

$ws2.Range("A3").Value = "BRVM - SERVICES FINANCIERS"
$ws2.Range("B3").Value = 3549.21

$ws2.Range("A4").Value = "BRVM-PRESTIGE"
$ws2.Range("B4").Value = 3377.7

$ws2.Range("A5").Value = "BRVM - INDUSTRIELS"
$ws2.Range("B5").Value = 2891.32

$ws2.Range("A6").Value = "BRVM - ENERGIE"
$ws2.Range("B6").Value = 1916.61

$ws2.Range("A7").Value = "BRVM - SERVICES PUBLICS"
$ws2.Range("B7").Value = 1739.57

$ws2.Range("A8").Value = "BRVM - TELECOMMUNICATIONS"
$ws2.Range("B8").Value = 1298.74

$ws2.Range("A9").Value = "EVIOSYS PACKAGING SIEM CI (SEMC)"
$ws2.Range("B9").Value = 31.68

$ws2.Range("A10").Value = "SICOR CI (SICC)"
$ws2.Range("B10").Value = 11.78

$ws2.Range("A11").Value = "NEI-CEDA CI (NEIC)"
$ws2.Range("B11").Value = 11.09
